$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with the default (unstyled) look, used to strip the text-number
# style flag Excel adds when a numeric-looking literal is forced to Text below.
$plainStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "26.704.64"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.597.63"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.24"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.822.29"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.635.29"
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.33"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").Value = "0.0₃0767"
$ws.Range("E17").Value = "  +5.48%  "
$ws.Range("D18").Value = "26.671.13"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "209.37"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.11"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +4.75%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.31"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.94"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.11"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.14"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.32"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("D34").Value = "1.284.50"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  -6.28%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +16.05%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.44"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.783"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.26"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "1.734.90"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.34"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.100"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.34"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -1.35%  "
